$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Date column (BF) keeps text formatting so the corrected
# "YYYY-MM-DD" strings are not auto-converted into date serial numbers.
$ws.Range("BF2:BF31").NumberFormat = "@"

# Per-row numeric stat corrections plus the BF (Date) text fix, derived
# from the authoritative box-score re-pull (one day had been off by a day).

# Row 2
$ws.Range("D2").Value2 = 82
$ws.Range("E2").Value2 = 53
$ws.Range("G2").Value2 = 0.646
$ws.Range("I2").Value2 = 38.8
$ws.Range("K2").Value2 = 0.468
$ws.Range("M2").Value2 = 17.7
$ws.Range("N2").Value2 = 0.36
$ws.Range("R2").Value2 = 11.8
$ws.Range("T2").Value2 = 41.7
$ws.Range("U2").Value2 = 21.8
$ws.Range("W2").Value2 = 7.2
$ws.Range("Z2").Value2 = 19.9
$ws.Range("AA2").Value2 = 19.3
$ws.Range("AB2").Value2 = 101.7
$ws.Range("AC2").Value2 = 4.7
$ws.Range("AD2").Value2 = 1
$ws.Range("AE2").Value2 = 6
$ws.Range("AF2").Value2 = 6
$ws.Range("AG2").Value2 = 6
$ws.Range("AI2").Value2 = 6
$ws.Range("AN2").Value2 = 9
$ws.Range("AO2").Value2 = 21
$ws.Range("AR2").Value2 = 6
$ws.Range("AS2").Value2 = 23
$ws.Range("AT2").Value2 = 17
$ws.Range("AW2").Value2 = 15
$ws.Range("AZ2").Value2 = 8
$ws.Range("BB2").Value2 = 13
$ws.Range("BF2").Value2 = "2010-05-07"

# Row 3
$ws.Range("D3").Value2 = 82
$ws.Range("E3").Value2 = 50
$ws.Range("G3").Value2 = 0.61
$ws.Range("H3").Value2 = 48.2
$ws.Range("I3").Value2 = 37.1
$ws.Range("K3").Value2 = 0.483
$ws.Range("L3").Value2 = 6.1
$ws.Range("M3").Value2 = 17.5
$ws.Range("N3").Value2 = 0.348
$ws.Range("P3").Value2 = 25.5
$ws.Range("Q3").Value2 = 0.746
$ws.Range("R3").Value2 = 8.699999999999999
$ws.Range("S3").Value2 = 29.9
$ws.Range("T3").Value2 = 38.6
$ws.Range("X3").Value2 = 4.9
$ws.Range("AB3").Value2 = 99.2
$ws.Range("AC3").Value2 = 3.7
$ws.Range("AD3").Value2 = 1
$ws.Range("AE3").Value2 = 9
$ws.Range("AG3").Value2 = 9
$ws.Range("AL3").Value2 = 16
$ws.Range("AM3").Value2 = 16
$ws.Range("AN3").Value2 = 17
$ws.Range("AO3").Value2 = 15
$ws.Range("AQ3").Value2 = 21
$ws.Range("AX3").Value2 = 15
$ws.Range("AY3").Value2 = 15
$ws.Range("AZ3").Value2 = 22
$ws.Range("BC3").Value2 = 9
$ws.Range("BF3").Value2 = "2010-05-07"

# Row 4
$ws.Range("D4").Value2 = 82
$ws.Range("F4").Value2 = 38
$ws.Range("G4").Value2 = 0.537
$ws.Range("I4").Value2 = 34.9
$ws.Range("K4").Value2 = 0.453
$ws.Range("L4").Value2 = 5.6
$ws.Range("M4").Value2 = 16.2
$ws.Range("N4").Value2 = 0.346
$ws.Range("Q4").Value2 = 0.751
$ws.Range("R4").Value2 = 10.5
$ws.Range("S4").Value2 = 30.4
$ws.Range("U4").Value2 = 20.2
$ws.Range("X4").Value2 = 5.4
$ws.Range("Z4").Value2 = 19.5
$ws.Range("AB4").Value2 = 95.3
$ws.Range("AC4").Value2 = 1.5
$ws.Range("AD4").Value2 = 1
$ws.Range("AK4").Value2 = 22
$ws.Range("AN4").Value2 = 21
$ws.Range("AP4").Value2 = 5
$ws.Range("AQ4").Value2 = 20
$ws.Range("AR4").Value2 = 21
$ws.Range("AX4").Value2 = 7
$ws.Range("AY4").Value2 = 29
$ws.Range("AZ4").Value2 = 5
$ws.Range("BB4").Value2 = 28
$ws.Range("BC4").Value2 = 15
$ws.Range("BF4").Value2 = "2010-05-07"

# Row 5
$ws.Range("AE5").Value2 = 17
$ws.Range("AH5").Value2 = 2
$ws.Range("AJ5").Value2 = 10
$ws.Range("AV5").Value2 = 15
$ws.Range("AW5").Value2 = 24
$ws.Range("AY5").Value2 = 23
$ws.Range("BF5").Value2 = "2010-05-07"

# Row 6
$ws.Range("D6").Value2 = 82
$ws.Range("F6").Value2 = 21
$ws.Range("G6").Value2 = 0.744
$ws.Range("H6").Value2 = 48.2
$ws.Range("I6").Value2 = 37.8
$ws.Range("J6").Value2 = 77.90000000000001
$ws.Range("K6").Value2 = 0.485
$ws.Range("N6").Value2 = 0.381
$ws.Range("O6").Value2 = 19.1
$ws.Range("P6").Value2 = 26.6
$ws.Range("R6").Value2 = 9.6
$ws.Range("S6").Value2 = 32.8
$ws.Range("U6").Value2 = 22.4
$ws.Range("Y6").Value2 = 4
$ws.Range("Z6").Value2 = 19.4
$ws.Range("AA6").Value2 = 20.8
$ws.Range("AB6").Value2 = 102.1
$ws.Range("AC6").Value2 = 6.5
$ws.Range("AD6").Value2 = 1
$ws.Range("AI6").Value2 = 15
$ws.Range("AJ6").Value2 = 28
$ws.Range("AO6").Value2 = 12
$ws.Range("AP6").Value2 = 5
$ws.Range("AU6").Value2 = 6
$ws.Range("AV6").Value2 = 12
$ws.Range("BA6").Value2 = 15
$ws.Range("BF6").Value2 = "2010-05-07"

# Row 7
$ws.Range("D7").Value2 = 82
$ws.Range("F7").Value2 = 27
$ws.Range("G7").Value2 = 0.671
$ws.Range("I7").Value2 = 38.3
$ws.Range("J7").Value2 = 82.40000000000001
$ws.Range("K7").Value2 = 0.464
$ws.Range("N7").Value2 = 0.372
$ws.Range("O7").Value2 = 18.6
$ws.Range("P7").Value2 = 22.8
$ws.Range("U7").Value2 = 23.4
$ws.Range("V7").Value2 = 12.9
$ws.Range("X7").Value2 = 5.5
$ws.Range("Z7").Value2 = 19.1
$ws.Range("AA7").Value2 = 20.1
$ws.Range("AB7").Value2 = 102
$ws.Range("AC7").Value2 = 2.7
$ws.Range("AD7").Value2 = 1
$ws.Range("AI7").Value2 = 11
$ws.Range("AN7").Value2 = 5
$ws.Range("AP7").Value2 = 25
$ws.Range("AX7").Value2 = 6
$ws.Range("BA7").Value2 = 21
$ws.Range("BF7").Value2 = "2010-05-07"

# Row 8
$ws.Range("D8").Value2 = 82
$ws.Range("E8").Value2 = 53
$ws.Range("G8").Value2 = 0.646
$ws.Range("I8").Value2 = 38.1
$ws.Range("J8").Value2 = 81.40000000000001
$ws.Range("L8").Value2 = 6.6
$ws.Range("M8").Value2 = 18.5
$ws.Range("N8").Value2 = 0.359
$ws.Range("S8").Value2 = 30.5
$ws.Range("W8").Value2 = 8.300000000000001
$ws.Range("Y8").Value2 = 5.3
$ws.Range("AB8").Value2 = 106.5
$ws.Range("AC8").Value2 = 4.1
$ws.Range("AD8").Value2 = 1
$ws.Range("AE8").Value2 = 6
$ws.Range("AF8").Value2 = 6
$ws.Range("AG8").Value2 = 6
$ws.Range("AH8").Value2 = 17
$ws.Range("AI8").Value2 = 13
$ws.Range("AN8").Value2 = 10
$ws.Range("AR8").Value2 = 17
$ws.Range("AS8").Value2 = 16
$ws.Range("AV8").Value2 = 11
$ws.Range("AY8").Value2 = 23
$ws.Range("BF8").Value2 = "2010-05-07"

# Row 9
$ws.Range("D9").Value2 = 82
$ws.Range("E9").Value2 = 27
$ws.Range("G9").Value2 = 0.329
$ws.Range("I9").Value2 = 35.9
$ws.Range("J9").Value2 = 80.5
$ws.Range("L9").Value2 = 4.6
$ws.Range("N9").Value2 = 0.314
$ws.Range("O9").Value2 = 17.7
$ws.Range("P9").Value2 = 24.4
$ws.Range("Q9").Value2 = 0.728
$ws.Range("R9").Value2 = 12.8
$ws.Range("V9").Value2 = 13.4
$ws.Range("AA9").Value2 = 20.8
$ws.Range("AC9").Value2 = -5.1
$ws.Range("AD9").Value2 = 1
$ws.Range("AE9").Value2 = 24
$ws.Range("AF9").Value2 = 24
$ws.Range("AG9").Value2 = 24
$ws.Range("AH9").Value2 = 17
$ws.Range("AM9").Value2 = 26
$ws.Range("AT9").Value2 = 27
$ws.Range("AV9").Value2 = 7
$ws.Range("AW9").Value2 = 14
$ws.Range("BA9").Value2 = 16
$ws.Range("BF9").Value2 = "2010-05-07"

# Row 10
$ws.Range("D10").Value2 = 82
$ws.Range("F10").Value2 = 56
$ws.Range("G10").Value2 = 0.317
$ws.Range("L10").Value2 = 7.7
$ws.Range("M10").Value2 = 20.6
$ws.Range("O10").Value2 = 19.9
$ws.Range("P10").Value2 = 25.4
$ws.Range("Q10").Value2 = 0.782
$ws.Range("R10").Value2 = 9.199999999999999
$ws.Range("Y10").Value2 = 5
$ws.Range("AD10").Value2 = 1
$ws.Range("AE10").Value2 = 26
$ws.Range("AF10").Value2 = 26
$ws.Range("AG10").Value2 = 26
$ws.Range("AN10").Value2 = 4
$ws.Range("AQ10").Value2 = 5
$ws.Range("AU10").Value2 = 5
$ws.Range("BC10").Value2 = 22
$ws.Range("BF10").Value2 = "2010-05-07"

# Row 11
$ws.Range("D11").Value2 = 82
$ws.Range("E11").Value2 = 42
$ws.Range("F11").Value2 = 40
$ws.Range("G11").Value2 = 0.512
$ws.Range("H11").Value2 = 48.5
$ws.Range("I11").Value2 = 37.7
$ws.Range("J11").Value2 = 84.40000000000001
$ws.Range("K11").Value2 = 0.447
$ws.Range("R11").Value2 = 11.8
$ws.Range("S11").Value2 = 30.1
$ws.Range("T11").Value2 = 42
$ws.Range("V11").Value2 = 14.5
$ws.Range("W11").Value2 = 7.1
$ws.Range("Y11").Value2 = 6.5
$ws.Range("Z11").Value2 = 20.9
$ws.Range("AB11").Value2 = 102.4
$ws.Range("AC11").Value2 = -0.4
$ws.Range("AD11").Value2 = 1
$ws.Range("AL11").Value2 = 5
$ws.Range("AO11").Value2 = 14
$ws.Range("AS11").Value2 = 21
$ws.Range("AW11").Value2 = 18
$ws.Range("AZ11").Value2 = 17
$ws.Range("BF11").Value2 = "2010-05-07"

# Row 12
$ws.Range("D12").Value2 = 82
$ws.Range("F12").Value2 = 50
$ws.Range("G12").Value2 = 0.39
$ws.Range("J12").Value2 = 83.2
$ws.Range("K12").Value2 = 0.443
$ws.Range("M12").Value2 = 23.1
$ws.Range("N12").Value2 = 0.348
$ws.Range("O12").Value2 = 19.1
$ws.Range("Q12").Value2 = 0.775
$ws.Range("V12").Value2 = 15
$ws.Range("W12").Value2 = 7.1
$ws.Range("AB12").Value2 = 100.8
$ws.Range("AC12").Value2 = -3
$ws.Range("AD12").Value2 = 1
$ws.Range("AH12").Value2 = 29
$ws.Range("AJ12").Value2 = 11
$ws.Range("AN12").Value2 = 18
$ws.Range("AO12").Value2 = 13
$ws.Range("AS12").Value2 = 6
$ws.Range("AU12").Value2 = 15
$ws.Range("AV12").Value2 = 25
$ws.Range("AW12").Value2 = 17
$ws.Range("AY12").Value2 = 18
$ws.Range("BF12").Value2 = "2010-05-07"

# Row 13
$ws.Range("D13").Value2 = 82
$ws.Range("F13").Value2 = 53
$ws.Range("G13").Value2 = 0.354
$ws.Range("I13").Value2 = 36.6
$ws.Range("K13").Value2 = 0.455
$ws.Range("M13").Value2 = 17.8
$ws.Range("O13").Value2 = 16.6
$ws.Range("Q13").Value2 = 0.73
$ws.Range("U13").Value2 = 22.1
$ws.Range("V13").Value2 = 15.7
$ws.Range("W13").Value2 = 6.5
$ws.Range("Z13").Value2 = 19.3
$ws.Range("AB13").Value2 = 95.7
$ws.Range("AD13").Value2 = 1
$ws.Range("AK13").Value2 = 20
$ws.Range("AO13").Value2 = 28
$ws.Range("AP13").Value2 = 26
$ws.Range("AS13").Value2 = 17
$ws.Range("AT13").Value2 = 14
$ws.Range("AU13").Value2 = 9
$ws.Range("AW13").Value2 = 23
$ws.Range("AY13").Value2 = 9
$ws.Range("BB13").Value2 = 27
$ws.Range("BF13").Value2 = "2010-05-07"

# Row 14
$ws.Range("D14").Value2 = 82
$ws.Range("E14").Value2 = 57
$ws.Range("G14").Value2 = 0.695
$ws.Range("I14").Value2 = 38.3
$ws.Range("K14").Value2 = 0.457
$ws.Range("M14").Value2 = 19
$ws.Range("N14").Value2 = 0.341
$ws.Range("O14").Value2 = 18.5
$ws.Range("P14").Value2 = 24.2
$ws.Range("Q14").Value2 = 0.765
$ws.Range("R14").Value2 = 11.9
$ws.Range("V14").Value2 = 13.4
$ws.Range("W14").Value2 = 7.5
$ws.Range("Z14").Value2 = 19.4
$ws.Range("AA14").Value2 = 21.2
$ws.Range("AD14").Value2 = 1
$ws.Range("AI14").Value2 = 9
$ws.Range("AK14").Value2 = 18
$ws.Range("AP14").Value2 = 17
$ws.Range("AR14").Value2 = 4
$ws.Range("AX14").Value2 = 16
$ws.Range("AY14").Value2 = 9
$ws.Range("AZ14").Value2 = 4
$ws.Range("BF14").Value2 = "2010-05-07"

# Row 15
$ws.Range("D15").Value2 = 82
$ws.Range("F15").Value2 = 42
$ws.Range("G15").Value2 = 0.488
$ws.Range("H15").Value2 = 48.5
$ws.Range("I15").Value2 = 39.3
$ws.Range("J15").Value2 = 83.8
$ws.Range("K15").Value2 = 0.469
$ws.Range("N15").Value2 = 0.337
$ws.Range("P15").Value2 = 26.9
$ws.Range("Q15").Value2 = 0.733
$ws.Range("T15").Value2 = 43.5
$ws.Range("U15").Value2 = 18.8
$ws.Range("Z15").Value2 = 20.2
$ws.Range("AB15").Value2 = 102.5
$ws.Range("AC15").Value2 = -1.5
$ws.Range("AD15").Value2 = 1
$ws.Range("AF15").Value2 = 18
$ws.Range("AH15").Value2 = 2
$ws.Range("AI15").Value2 = 4
$ws.Range("AJ15").Value2 = 7
$ws.Range("AK15").Value2 = 9
$ws.Range("AN15").Value2 = 26
$ws.Range("AO15").Value2 = 8
$ws.Range("AU15").Value2 = 29
$ws.Range("AX15").Value2 = 18
$ws.Range("AY15").Value2 = 28
$ws.Range("BA15").Value2 = 2
$ws.Range("BB15").Value2 = 7
$ws.Range("BF15").Value2 = "2010-05-07"

# Row 16
$ws.Range("D16").Value2 = 82
$ws.Range("E16").Value2 = 47
$ws.Range("G16").Value2 = 0.573
$ws.Range("I16").Value2 = 36.4
$ws.Range("J16").Value2 = 79.5
$ws.Range("K16").Value2 = 0.458
$ws.Range("N16").Value2 = 0.346
$ws.Range("O16").Value2 = 17.7
$ws.Range("P16").Value2 = 23.5
$ws.Range("Q16").Value2 = 0.752
$ws.Range("V16").Value2 = 13.2
$ws.Range("AB16").Value2 = 96.5
$ws.Range("AC16").Value2 = 2.3
$ws.Range("AD16").Value2 = 1
$ws.Range("AI16").Value2 = 26
$ws.Range("AK16").Value2 = 17
$ws.Range("AL16").Value2 = 17
$ws.Range("AM16").Value2 = 17
$ws.Range("AN16").Value2 = 19
$ws.Range("AQ16").Value2 = 19
$ws.Range("AT16").Value2 = 15
$ws.Range("AU16").Value2 = 28
$ws.Range("AW16").Value2 = 12
$ws.Range("AZ16").Value2 = 16
$ws.Range("BF16").Value2 = "2010-05-07"

# Row 17
$ws.Range("AE17").Value2 = 14
$ws.Range("AL17").Value2 = 6
$ws.Range("AW17").Value2 = 19
$ws.Range("AX17").Value2 = 20
$ws.Range("BB17").Value2 = 23
$ws.Range("BC17").Value2 = 14
$ws.Range("BF17").Value2 = "2010-05-07"

# Row 18
$ws.Range("D18").Value2 = 82
$ws.Range("E18").Value2 = 15
$ws.Range("G18").Value2 = 0.183
$ws.Range("J18").Value2 = 84.40000000000001
$ws.Range("K18").Value2 = 0.449
$ws.Range("L18").Value2 = 4.9
$ws.Range("M18").Value2 = 14.4
$ws.Range("N18").Value2 = 0.341
$ws.Range("O18").Value2 = 17.5
$ws.Range("P18").Value2 = 23.5
$ws.Range("T18").Value2 = 42.9
$ws.Range("U18").Value2 = 19.8
$ws.Range("W18").Value2 = 7.3
$ws.Range("X18").Value2 = 3.7
$ws.Range("Y18").Value2 = 5.4
$ws.Range("AA18").Value2 = 20.6
$ws.Range("AC18").Value2 = -9.6
$ws.Range("AD18").Value2 = 1
$ws.Range("AH18").Value2 = 17
$ws.Range("AI18").Value2 = 14
$ws.Range("AJ18").Value2 = 3
$ws.Range("AK18").Value2 = 25
$ws.Range("AM18").Value2 = 28
$ws.Range("AQ18").Value2 = 22
$ws.Range("AR18").Value2 = 9
$ws.Range("AW18").Value2 = 13
$ws.Range("AY18").Value2 = 25
$ws.Range("BA18").Value2 = 17
$ws.Range("BF18").Value2 = "2010-05-07"

# Row 19
$ws.Range("D19").Value2 = 82
$ws.Range("F19").Value2 = 70
$ws.Range("G19").Value2 = 0.146
$ws.Range("L19").Value2 = 4.6
$ws.Range("N19").Value2 = 0.318
$ws.Range("O19").Value2 = 19.2
$ws.Range("P19").Value2 = 24.6
$ws.Range("Q19").Value2 = 0.78
$ws.Range("S19").Value2 = 28.8
$ws.Range("T19").Value2 = 39.7
$ws.Range("V19").Value2 = 14.4
$ws.Range("X19").Value2 = 4.8
$ws.Range("Y19").Value2 = 5.1
$ws.Range("AA19").Value2 = 20.1
$ws.Range("AC19").Value2 = -9.1
$ws.Range("AD19").Value2 = 1
$ws.Range("AH19").Value2 = 17
$ws.Range("AM19").Value2 = 27
$ws.Range("AQ19").Value2 = 6
$ws.Range("AV19").Value2 = 16
$ws.Range("AW19").Value2 = 20
$ws.Range("AY19").Value2 = 21
$ws.Range("AZ19").Value2 = 10
$ws.Range("BA19").Value2 = 23
$ws.Range("BF19").Value2 = "2010-05-07"

# Row 20
$ws.Range("D20").Value2 = 82
$ws.Range("F20").Value2 = 45
$ws.Range("G20").Value2 = 0.451
$ws.Range("J20").Value2 = 83.40000000000001
$ws.Range("M20").Value2 = 19.2
$ws.Range("N20").Value2 = 0.363
$ws.Range("P20").Value2 = 20.3
$ws.Range("Q20").Value2 = 0.778
$ws.Range("R20").Value2 = 10.4
$ws.Range("T20").Value2 = 40.3
$ws.Range("U20").Value2 = 22.3
$ws.Range("W20").Value2 = 7.6
$ws.Range("AA20").Value2 = 19.5
$ws.Range("AB20").Value2 = 100.2
$ws.Range("AC20").Value2 = -2.5
$ws.Range("AD20").Value2 = 1
$ws.Range("AI20").Value2 = 7
$ws.Range("AN20").Value2 = 8
$ws.Range("AQ20").Value2 = 7
$ws.Range("AR20").Value2 = 22
$ws.Range("AS20").Value2 = 24
$ws.Range("AT20").Value2 = 25
$ws.Range("AU20").Value2 = 8
$ws.Range("AY20").Value2 = 12
$ws.Range("BF20").Value2 = "2010-05-07"

# Row 21
$ws.Range("D21").Value2 = 82
$ws.Range("F21").Value2 = 53
$ws.Range("G21").Value2 = 0.354
$ws.Range("I21").Value2 = 38.1
$ws.Range("J21").Value2 = 83.90000000000001
$ws.Range("K21").Value2 = 0.455
$ws.Range("L21").Value2 = 9.1
$ws.Range("M21").Value2 = 26.2
$ws.Range("N21").Value2 = 0.346
$ws.Range("AB21").Value2 = 102.1
$ws.Range("AC21").Value2 = -3.8
$ws.Range("AD21").Value2 = 1
$ws.Range("AF21").Value2 = 22
$ws.Range("AG21").Value2 = 22
$ws.Range("AJ21").Value2 = 6
$ws.Range("AK21").Value2 = 21
$ws.Range("AN21").Value2 = 20
$ws.Range("AQ21").Value2 = 4
$ws.Range("AS21").Value2 = 20
$ws.Range("AY21").Value2 = 13
$ws.Range("AZ21").Value2 = 9
$ws.Range("BA21").Value2 = 29
$ws.Range("BB21").Value2 = 9
$ws.Range("BC21").Value2 = 23
$ws.Range("BF21").Value2 = "2010-05-07"

# Row 22
$ws.Range("D22").Value2 = 82
$ws.Range("E22").Value2 = 50
$ws.Range("G22").Value2 = 0.61
$ws.Range("I22").Value2 = 37.4
$ws.Range("J22").Value2 = 80.8
$ws.Range("K22").Value2 = 0.462
$ws.Range("N22").Value2 = 0.34
$ws.Range("O22").Value2 = 21.7
$ws.Range("R22").Value2 = 11.7
$ws.Range("AC22").Value2 = 3.5
$ws.Range("AD22").Value2 = 1
$ws.Range("AN22").Value2 = 25
$ws.Range("AR22").Value2 = 10
$ws.Range("AY22").Value2 = 14
$ws.Range("AZ22").Value2 = 20
$ws.Range("BB22").Value2 = 14
$ws.Range("BC22").Value2 = 10
$ws.Range("BF22").Value2 = "2010-05-07"

# Row 23
$ws.Range("D23").Value2 = 82
$ws.Range("E23").Value2 = 59
$ws.Range("G23").Value2 = 0.72
$ws.Range("I23").Value2 = 36.6
$ws.Range("J23").Value2 = 78
$ws.Range("K23").Value2 = 0.47
$ws.Range("L23").Value2 = 10.3
$ws.Range("N23").Value2 = 0.375
$ws.Range("O23").Value2 = 19.2
$ws.Range("P23").Value2 = 26.5
$ws.Range("Q23").Value2 = 0.724
$ws.Range("T23").Value2 = 43.2
$ws.Range("U23").Value2 = 19.7
$ws.Range("V23").Value2 = 14.1
$ws.Range("W23").Value2 = 6.2
$ws.Range("AB23").Value2 = 102.8
$ws.Range("AC23").Value2 = 7.5
$ws.Range("AD23").Value2 = 1
$ws.Range("AJ23").Value2 = 27
$ws.Range("AK23").Value2 = 7
$ws.Range("AN23").Value2 = 3
$ws.Range("AP23").Value2 = 7
$ws.Range("AZ23").Value2 = 7
$ws.Range("BB23").Value2 = 6
$ws.Range("BF23").Value2 = "2010-05-07"

# Row 24
$ws.Range("D24").Value2 = 82
$ws.Range("F24").Value2 = 55
$ws.Range("G24").Value2 = 0.329
$ws.Range("O24").Value2 = 16.7
$ws.Range("S24").Value2 = 29.5
$ws.Range("T24").Value2 = 41
$ws.Range("V24").Value2 = 14.5
$ws.Range("X24").Value2 = 5.4
$ws.Range("Z24").Value2 = 20.5
$ws.Range("AB24").Value2 = 97.7
$ws.Range("AC24").Value2 = -3.9
$ws.Range("AD24").Value2 = 1
$ws.Range("AO24").Value2 = 27
$ws.Range("AX24").Value2 = 8
$ws.Range("BB24").Value2 = 22
$ws.Range("BF24").Value2 = "2010-05-07"

# Row 25
$ws.Range("D25").Value2 = 82
$ws.Range("E25").Value2 = 54
$ws.Range("G25").Value2 = 0.659
$ws.Range("M25").Value2 = 21.6
$ws.Range("Q25").Value2 = 0.77
$ws.Range("S25").Value2 = 31.9
$ws.Range("T25").Value2 = 43
$ws.Range("U25").Value2 = 23.3
$ws.Range("Y25").Value2 = 4.5
$ws.Range("AD25").Value2 = 1
$ws.Range("AO25").Value2 = 5
$ws.Range("AP25").Value2 = 9
$ws.Range("AS25").Value2 = 7
$ws.Range("AY25").Value2 = 11
$ws.Range("AZ25").Value2 = 18
$ws.Range("BC25").Value2 = 5
$ws.Range("BF25").Value2 = "2010-05-07"

# Row 26
$ws.Range("D26").Value2 = 82
$ws.Range("E26").Value2 = 50
$ws.Range("G26").Value2 = 0.61
$ws.Range("L26").Value2 = 6
$ws.Range("N26").Value2 = 0.354
$ws.Range("Q26").Value2 = 0.79
$ws.Range("S26").Value2 = 29.1
$ws.Range("T26").Value2 = 40.2
$ws.Range("V26").Value2 = 12.3
$ws.Range("W26").Value2 = 6.4
$ws.Range("X26").Value2 = 4.3
$ws.Range("Z26").Value2 = 20.9
$ws.Range("AC26").Value2 = 3.3
$ws.Range("AD26").Value2 = 1
$ws.Range("AI26").Value2 = 27
$ws.Range("AN26").Value2 = 13
$ws.Range("AT26").Value2 = 26
$ws.Range("AZ26").Value2 = 19
$ws.Range("BF26").Value2 = "2010-05-07"

# Row 27
$ws.Range("D27").Value2 = 82
$ws.Range("F27").Value2 = 57
$ws.Range("G27").Value2 = 0.305
$ws.Range("I27").Value2 = 38.3
$ws.Range("K27").Value2 = 0.456
$ws.Range("N27").Value2 = 0.349
$ws.Range("Q27").Value2 = 0.726
$ws.Range("S27").Value2 = 30.7
$ws.Range("T27").Value2 = 42.6
$ws.Range("U27").Value2 = 20.5
$ws.Range("AB27").Value2 = 100
$ws.Range("AC27").Value2 = -4.4
$ws.Range("AD27").Value2 = 1
$ws.Range("AE27").Value2 = 28
$ws.Range("AF27").Value2 = 28
$ws.Range("AG27").Value2 = 28
$ws.Range("AI27").Value2 = 9
$ws.Range("AK27").Value2 = 19
$ws.Range("AM27").Value2 = 20
$ws.Range("AP27").Value2 = 18
$ws.Range("AV27").Value2 = 23
$ws.Range("AW27").Value2 = 21
$ws.Range("BA27").Value2 = 22
$ws.Range("BF27").Value2 = "2010-05-07"

# Row 28
$ws.Range("D28").Value2 = 82
$ws.Range("E28").Value2 = 50
$ws.Range("G28").Value2 = 0.61
$ws.Range("I28").Value2 = 38.4
$ws.Range("J28").Value2 = 81.2
$ws.Range("K28").Value2 = 0.473
$ws.Range("L28").Value2 = 6.8
$ws.Range("M28").Value2 = 18.9
$ws.Range("P28").Value2 = 24
$ws.Range("Q28").Value2 = 0.74
$ws.Range("U28").Value2 = 22.3
$ws.Range("V28").Value2 = 13.6
$ws.Range("X28").Value2 = 4.6
$ws.Range("Z28").Value2 = 20.4
$ws.Range("AB28").Value2 = 101.4
$ws.Range("AC28").Value2 = 5.1
$ws.Range("AD28").Value2 = 1
$ws.Range("AH28").Value2 = 17
$ws.Range("AI28").Value2 = 8
$ws.Range("AP28").Value2 = 18
$ws.Range("AQ28").Value2 = 24
$ws.Range("AR28").Value2 = 18
$ws.Range("AU28").Value2 = 7
$ws.Range("AX28").Value2 = 22
$ws.Range("BF28").Value2 = "2010-05-07"

# Row 29
$ws.Range("D29").Value2 = 82
$ws.Range("E29").Value2 = 40
$ws.Range("G29").Value2 = 0.488
$ws.Range("I29").Value2 = 39
$ws.Range("K29").Value2 = 0.482
$ws.Range("O29").Value2 = 19.7
$ws.Range("P29").Value2 = 25.8
$ws.Range("R29").Value2 = 9.800000000000001
$ws.Range("S29").Value2 = 30.6
$ws.Range("T29").Value2 = 40.4
$ws.Range("Z29").Value2 = 22.2
$ws.Range("AA29").Value2 = 21.1
$ws.Range("AC29").Value2 = -1.8
$ws.Range("AD29").Value2 = 1
$ws.Range("AE29").Value2 = 18
$ws.Range("AF29").Value2 = 18
$ws.Range("AG29").Value2 = 18
$ws.Range("AH29").Value2 = 17
$ws.Range("AO29").Value2 = 7
$ws.Range("AP29").Value2 = 8
$ws.Range("AS29").Value2 = 15
$ws.Range("AT29").Value2 = 23
$ws.Range("AU29").Value2 = 10
$ws.Range("AX29").Value2 = 21
$ws.Range("AZ29").Value2 = 23
$ws.Range("BF29").Value2 = "2010-05-07"

# Row 30
$ws.Range("D30").Value2 = 82
$ws.Range("F30").Value2 = 29
$ws.Range("G30").Value2 = 0.646
$ws.Range("L30").Value2 = 5.4
$ws.Range("M30").Value2 = 14.7
$ws.Range("N30").Value2 = 0.364
$ws.Range("O30").Value2 = 20.2
$ws.Range("Q30").Value2 = 0.741
$ws.Range("W30").Value2 = 8.199999999999999
$ws.Range("AA30").Value2 = 22.2
$ws.Range("AC30").Value2 = 5.3
$ws.Range("AD30").Value2 = 1
$ws.Range("AE30").Value2 = 6
$ws.Range("AF30").Value2 = 6
$ws.Range("AG30").Value2 = 6
$ws.Range("AI30").Value2 = 3
$ws.Range("AN30").Value2 = 7
$ws.Range("AQ30").Value2 = 23
$ws.Range("AX30").Value2 = 16
$ws.Range("BF30").Value2 = "2010-05-07"

# Row 31
$ws.Range("D31").Value2 = 82
$ws.Range("E31").Value2 = 26
$ws.Range("G31").Value2 = 0.317
$ws.Range("M31").Value2 = 14.9
$ws.Range("Q31").Value2 = 0.762
$ws.Range("S31").Value2 = 30
$ws.Range("T31").Value2 = 41.8
$ws.Range("Y31").Value2 = 5.1
$ws.Range("Z31").Value2 = 21.4
$ws.Range("AC31").Value2 = -4.8
$ws.Range("AD31").Value2 = 1
$ws.Range("AE31").Value2 = 26
$ws.Range("AF31").Value2 = 26
$ws.Range("AG31").Value2 = 26
$ws.Range("AK31").Value2 = 24
$ws.Range("AN31").Value2 = 14
$ws.Range("AR31").Value2 = 7
$ws.Range("AS31").Value2 = 22
$ws.Range("AT31").Value2 = 16
$ws.Range("AV31").Value2 = 21
$ws.Range("AY31").Value2 = 19
$ws.Range("AZ31").Value2 = 21
$ws.Range("BA31").Value2 = 18
$ws.Range("BF31").Value2 = "2010-05-07"
